$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 82 is no longer the most recent day, so it loses the "last row" date-only
# style and reverts to the regular date+time style used by all prior rows.
$ws.Range("A82").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's row (daily update).
$ws.Range("A83").Value = 45670
$ws.Range("A83").NumberFormat = "YYYY-MM-DD"

$ws.Range("B83").Value = 195
$ws.Range("C83").Value = 193
$ws.Range("D83").Value = 191
